# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly computed "K" (strikeout) values replacing the old Strike# derived values
$kVals = @(1, 0, 0, 1, 1, 1, 3, 1, 1, 1, 0, 0, 0, 2, 0, 2, 2, 1, 1, 2, 1)

$row = 2
foreach ($k in $kVals) {
    $ws.Cells.Item($row, 7).Value = $k
    $row++
}
